# Remove zero padding from the hour portion of the "HH:MM-HH:MM" time-range
# strings used throughout the schedule workbook (e.g. "02:00-02:08" ->
# "2:00-2:08"). The same text values are duplicated (via the shared string
# table) across the Match Schedule, Judging Schedule and every individual
# Team Schedule sheet, so every worksheet is scanned and updated.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
            $val = $cell.Value2

            if ($val -is [string] -and $val -match '^0?\d:\d{2}-0?\d:\d{2}$') {
                # Strip a leading zero from each "0H:" hour segment.
                $newVal = $val -replace '(^|-)0(\d:)', '${1}${2}'
                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
